$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - column headers (labels reshuffled/renamed)
$ws.Range("A1").Value = "Comarca nombre"
$ws.Range("B1").Value = "Número hogares"
$ws.Range("C1").Value = "Comarca código"
$ws.Range("D1").Value = "Provincia código"
$ws.Range("E1").Value = "Aragón"
$ws.Range("F1").Value = "Núcleos en el hogar"
$ws.Range("G1").Value = "Municipio código"
$ws.Range("H1").Value = "Provincia nombre"
$ws.Range("I1").Value = "Municipio nombre"

# Row 2 - "concept" row (sdmx-dimension / iaest-measure / null)
$ws.Range("A2").Value = "sdmx-dimension:refArea"
$ws.Range("B2").Value = "iaest-measure:numero-hogares"
$ws.Range("D2").Value = "null"
$ws.Range("F2").Value = "iaest-measure:nucleos-en-el-hogar"
$ws.Range("I2").Value = "sdmx-dimension:refArea"

# Row 3 - "role" row (dim / medida / null)
$ws.Range("A3").Value = "dim"
$ws.Range("D3").Value = "null"
$ws.Range("F3").Value = "medida"
$ws.Range("I3").Value = "dim"

# Row 4 - "type/URI" row
$ws.Range("A4").Value = "URI-comarca"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("F4").Value = "xsd:string"
$ws.Range("H4").Value = "URI-Provincia"
$ws.Range("I4").Value = "URI-Municipio"
